$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update resource request limits (Job Supervisor, Cron Scheduler, Socials)
$ws.Range("B5").Value = 0.2
$ws.Range("B6").Value = 0.2
$ws.Range("C6").Value = 0.5
$ws.Range("B7").Value = 0.2
$ws.Range("C7").Value = 0.5

# Move the active selection to C8
$ws.Range("C8").Select()
